$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "v3.85.0"
$ws.Range("E3").Value = "v2025-12-23"
$ws.Range("E4").Select()
